# Applies data-driven updates to the "Transferencia" workbook for the
# "interaction para validar carga de datos y clave dinámica" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# --- Row 3 (ID 2): change destination product type / account references ---
$ws.Range("N3").Value = "Corriente"
$ws.Range("Q3").Value = "406-182800-03"
$ws.Range("R3").Value = "Corriente"

# --- Row 7 (ID 6): update user, product and amount data ---
$ws.Range("D7").Value = "pruebauser01"
$ws.Range("E7").Value = "6789"
$ws.Range("M7").Value = "0437002003817"
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "FIDUCUENTA"
$ws.Range("O7").Value = 3000
$ws.Range("Q7").Value = "406-757180-07"
$ws.Range("R7").Value = "Ahorros"
$ws.Range("T7").Value = "Inscritos"

# --- Row 8 (ID 7): update destination product type and expected result ---
$ws.Range("N8").Value = "FIDUCUENTA"
$ws.Range("T8").Value = "Inscritos"

# --- Sheet view: scroll / selection moved from O7 to R8 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 16   # column P is left-most visible column
$ws.Range("R8").Select()
